# Auto-generated Excel COM-interop script
# Applies updated currentAveragePrice / LevePrice / LeveProfit figures
# (columns H:N) for specific Leve rows across all 8 job sheets,
# as produced by the scheduled market-data refresh runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 132.5
$ws.Range("I55").Value = 95
$ws.Range("J55").Value = 207.5
$ws.Range("K55").Value = 95
$ws.Range("L55").Value = 207.5
$ws.Range("M55").Value = 119
$ws.Range("N55").Value = -635.5
$ws.Range("H70").Value = 31432104
$ws.Range("I70").Value = 4547724.5
$ws.Range("J70").Value = 76928744
$ws.Range("K70").Value = 13643173.5
$ws.Range("L70").Value = 230786232
$ws.Range("M70").Value = -13642903.5
$ws.Range("N70").Value = -230786772
$ws.Range("H73").Value = 31432104
$ws.Range("I73").Value = 4547724.5
$ws.Range("J73").Value = 76928744
$ws.Range("K73").Value = 13643173.5
$ws.Range("L73").Value = 230786232
$ws.Range("M73").Value = -13642237.5
$ws.Range("N73").Value = -230788104
$ws.Range("H76").Value = 3334.818
$ws.Range("I76").Value = 3125.1875
$ws.Range("J76").Value = 3893.8333
$ws.Range("K76").Value = 3125.1875
$ws.Range("L76").Value = 3893.8333
$ws.Range("M76").Value = -2810.1875
$ws.Range("N76").Value = -4523.8333
$ws.Range("H79").Value = 3334.818
$ws.Range("I79").Value = 3125.1875
$ws.Range("J79").Value = 3893.8333
$ws.Range("K79").Value = 3125.1875
$ws.Range("L79").Value = 3893.8333
$ws.Range("M79").Value = -2033.1875
$ws.Range("N79").Value = -6077.8333
$ws.Range("H92").Value = 3636.0386
$ws.Range("I92").Value = 1536.7894
$ws.Range("J92").Value = 9334
$ws.Range("K92").Value = 1536.7894
$ws.Range("L92").Value = 9334
$ws.Range("M92").Value = -288.7893999999999
$ws.Range("N92").Value = -11830
$ws.Range("H132").Value = 3264.5881
$ws.Range("I132").Value = 3181.5454
$ws.Range("J132").Value = 3416.8333
$ws.Range("K132").Value = 9544.636200000001
$ws.Range("L132").Value = 10250.4999
$ws.Range("M132").Value = -7014.636200000001
$ws.Range("H138").Value = 2056.9795
$ws.Range("I138").Value = 1018.0476
$ws.Range("J138").Value = 2836.1785
$ws.Range("K138").Value = 3054.1428
$ws.Range("L138").Value = 8508.5355
$ws.Range("M138").Value = 2085.8572
$ws.Range("N138").Value = -18788.5355

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3769.4546
$ws.Range("I2").Value = 79.5
$ws.Range("J2").Value = 4589.4443
$ws.Range("K2").Value = 79.5
$ws.Range("L2").Value = 4589.4443
$ws.Range("M2").Value = 33.5
$ws.Range("N2").Value = -4815.4443
$ws.Range("H32").Value = 4767.1724
$ws.Range("I32").Value = 4990.815
$ws.Range("J32").Value = 1748
$ws.Range("K32").Value = 4990.815
$ws.Range("L32").Value = 1748
$ws.Range("M32").Value = -4703.815
$ws.Range("N32").Value = -2322
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").ClearContents()
$ws.Range("H50").Value = 2209.2
$ws.Range("I50").Value = 1224
$ws.Range("J50").Value = 2866
$ws.Range("K50").Value = 1224
$ws.Range("L50").Value = 2866
$ws.Range("M50").Value = -510
$ws.Range("N50").Value = -4294
$ws.Range("H61").Value = 5074.4346
$ws.Range("I61").Value = 4835.6
$ws.Range("J61").Value = 6666.6665
$ws.Range("K61").Value = 4835.6
$ws.Range("L61").Value = 6666.6665
$ws.Range("M61").Value = -4623.6
$ws.Range("N61").Value = -7090.6665
$ws.Range("H74").Value = 4059.0833
$ws.Range("I74").Value = 3899.125
$ws.Range("J74").Value = 4379
$ws.Range("K74").Value = 3899.125
$ws.Range("L74").Value = 4379
$ws.Range("M74").Value = -3025.125
$ws.Range("N74").Value = -6127
$ws.Range("H77").Value = 4059.0833
$ws.Range("I77").Value = 3899.125
$ws.Range("J77").Value = 4379
$ws.Range("K77").Value = 19495.625
$ws.Range("L77").Value = 21895
$ws.Range("M77").Value = -15127.625
$ws.Range("N77").Value = -30631
$ws.Range("H110").Value = 956.46155
$ws.Range("I110").Value = 936.1667
$ws.Range("J110").Value = 1200
$ws.Range("K110").Value = 936.1667
$ws.Range("L110").Value = 1200
$ws.Range("M110").Value = 1108.8333
$ws.Range("N110").Value = -5290
$ws.Range("H116").Value = 3769.4546
$ws.Range("I116").Value = 79.5
$ws.Range("J116").Value = 4589.4443
$ws.Range("K116").Value = 79.5
$ws.Range("L116").Value = 4589.4443
$ws.Range("M116").Value = 2214.5
$ws.Range("N116").Value = -9177.444299999999
$ws.Range("H132").Value = 2977.7932
$ws.Range("I132").Value = 2150.6365
$ws.Range("J132").Value = 5577.4287
$ws.Range("K132").Value = 6451.9095
$ws.Range("L132").Value = 16732.2861
$ws.Range("M132").Value = -3921.9095
$ws.Range("H136").Value = 5074.4346
$ws.Range("I136").Value = 4835.6
$ws.Range("J136").Value = 6666.6665
$ws.Range("K136").Value = 14506.8
$ws.Range("L136").Value = 19999.9995
$ws.Range("M136").Value = -11956.8
$ws.Range("N136").Value = -25099.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3769.4546
$ws.Range("I3").Value = 79.5
$ws.Range("J3").Value = 4589.4443
$ws.Range("K3").Value = 79.5
$ws.Range("L3").Value = 4589.4443
$ws.Range("M3").Value = 34.5
$ws.Range("N3").Value = -4817.4443
$ws.Range("H86").Value = 16130658
$ws.Range("I86").Value = 1694.6957
$ws.Range("J86").Value = 62501428
$ws.Range("K86").Value = 1694.6957
$ws.Range("L86").Value = 62501428
$ws.Range("M86").Value = -571.6957
$ws.Range("N86").Value = -62503674
$ws.Range("H89").Value = 16130658
$ws.Range("I89").Value = 1694.6957
$ws.Range("J89").Value = 62501428
$ws.Range("K89").Value = 8473.478499999999
$ws.Range("L89").Value = 312507140
$ws.Range("M89").Value = -2857.478499999999
$ws.Range("N89").Value = -312518372
$ws.Range("H94").Value = 4236.636
$ws.Range("I94").Value = 5001.8
$ws.Range("J94").Value = 3599
$ws.Range("K94").Value = 5001.8
$ws.Range("L94").Value = 3599
$ws.Range("M94").Value = -4550.8
$ws.Range("N94").Value = -4501
$ws.Range("H96").Value = 28472.428
$ws.Range("I96").Value = 9863.4
$ws.Range("J96").Value = 74995
$ws.Range("K96").Value = 9863.4
$ws.Range("L96").Value = 74995
$ws.Range("M96").Value = -7117.4
$ws.Range("N96").Value = -80487

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2915.2
$ws.Range("I16").Value = 2897.7144
$ws.Range("J16").Value = 2930.5
$ws.Range("K16").Value = 2897.7144
$ws.Range("L16").Value = 2930.5
$ws.Range("M16").Value = -2610.7144
$ws.Range("N16").Value = -3504.5
$ws.Range("H41").Value = 26219.777
$ws.Range("I41").Value = 10000
$ws.Range("J41").Value = 34329.668
$ws.Range("K41").Value = 10000
$ws.Range("L41").Value = 34329.668
$ws.Range("M41").Value = -9572
$ws.Range("N41").Value = -35185.668
$ws.Range("H113").Value = 2915.2
$ws.Range("I113").Value = 2897.7144
$ws.Range("J113").Value = 2930.5
$ws.Range("K113").Value = 2897.7144
$ws.Range("L113").Value = 2930.5
$ws.Range("M113").Value = -727.7143999999998
$ws.Range("N113").Value = -7270.5
$ws.Range("H134").Value = 2279.15
$ws.Range("I134").Value = 2105.1667
$ws.Range("J134").Value = 3845
$ws.Range("K134").Value = 6315.500100000001
$ws.Range("L134").Value = 11535
$ws.Range("M134").Value = -3780.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2636.3333
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 2636.3333
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 7908.999899999999
$ws.Range("N64").Value = -8448.999899999999
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 2636.3333
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 2636.3333
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 7908.999899999999
$ws.Range("N67").Value = -9780.999899999999
$ws.Range("M67").ClearContents()
$ws.Range("H87").Value = 7333.6665
$ws.Range("I87").Value = 5200.6
$ws.Range("J87").Value = 10000
$ws.Range("K87").Value = 15601.8
$ws.Range("L87").Value = 30000
$ws.Range("M87").Value = -14353.8
$ws.Range("N87").Value = -32496
$ws.Range("H90").Value = 7333.6665
$ws.Range("I90").Value = 5200.6
$ws.Range("J90").Value = 10000
$ws.Range("K90").Value = 46805.4
$ws.Range("L90").Value = 90000
$ws.Range("M90").Value = -40565.4
$ws.Range("N90").Value = -102480
$ws.Range("H113").Value = 1174.3572
$ws.Range("I113").Value = 700
$ws.Range("J113").Value = 1210.8462
$ws.Range("K113").Value = 2100
$ws.Range("L113").Value = 3632.5386
$ws.Range("M113").Value = 70
$ws.Range("N113").Value = -7972.5386
$ws.Range("H121").Value = 482.4
$ws.Range("I121").Value = 433.33334
$ws.Range("J121").Value = 522.5454999999999
$ws.Range("K121").Value = 1300.00002
$ws.Range("L121").Value = 1567.6365
$ws.Range("M121").Value = 9.99998000000005
$ws.Range("N121").Value = -4187.6365

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 75000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 75000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 75000
$ws.Range("N64").Value = -75496
$ws.Range("H67").Value = 75000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 75000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 75000
$ws.Range("N67").Value = -76716
$ws.Range("H132").Value = 2798.1538
$ws.Range("I132").Value = 2824.8262
$ws.Range("J132").Value = 2593.6667
$ws.Range("K132").Value = 8474.4786
$ws.Range("L132").Value = 7781.000100000001
$ws.Range("M132").Value = -5944.4786
$ws.Range("N132").Value = -12841.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 7856.857
$ws.Range("I93").Value = 7624
$ws.Range("J93").Value = 7950
$ws.Range("K93").Value = 7624
$ws.Range("L93").Value = 7950
$ws.Range("M93").Value = -6376
$ws.Range("N93").Value = -10446
$ws.Range("H132").Value = 3556.5134
$ws.Range("I132").Value = 3248.0967
$ws.Range("J132").Value = 5150
$ws.Range("K132").Value = 9744.2901
$ws.Range("L132").Value = 15450
$ws.Range("M132").Value = -7214.2901
$ws.Range("N132").Value = -20510

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 2000
$ws.Range("I54").Value = 2000
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 2000
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -1480
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H132").Value = 2464.5925
$ws.Range("I132").Value = 2517.652
$ws.Range("J132").Value = 2159.5
$ws.Range("K132").Value = 7552.956
$ws.Range("L132").Value = 6478.5
$ws.Range("M132").Value = -5022.956
